$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Plain value edits (style unchanged)
# ---------------------------------------------------------------------------
$ws.Range("E9").Value2 = "Ir a hacer mercado"
$ws.Range("E11").Value2 = "Cocinar"
$ws.Range("E15").Value2 = "GYM"
$ws.Range("C19").Value2 = "GYM"
$ws.Range("D21").Value2 = "Actividades Varias"

# ---------------------------------------------------------------------------
# Value + style edits (copy fill/border from a cell with the target style)
# ---------------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$ws.Range("G13").Value2 = "Actividades Varias"

$ws.Range("B3").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("G14").Value2 = "Actividades Varias"

$ws.Range("I7").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value2 = "Taller de fisica mecanica"

$ws.Range("D13").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value2 = "Estudiar Algebra lineal"

$ws.Range("I7").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value2 = "Taller de fisica mecanica"

$ws.Range("B3").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value2 = "GYM"

$ws.Range("B18").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value2 = "Cenar"

$ws.Range("B3").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value2 = "Actividades Varias"

$ws.Range("B3").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D19").Value2 = "GYM"

$ws.Range("B3").Copy()
$ws.Range("G19").PasteSpecial(-4122)
$ws.Range("G19").Value2 = "Actividades Varias"

$ws.Range("B3").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value2 = "Actividades Varias"

$ws.Range("B3").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value2 = "Actividades Varias"

$ws.Range("D13").Copy()
$ws.Range("G20").PasteSpecial(-4122)
$ws.Range("G20").Value2 = "Estudiar Algebra lineal"

$ws.Range("D9").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("C21").Value2 = "Ingles "

$ws.Range("D9").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value2 = "Ingles "

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Selection change: active cell moves from B9 to G7
# ---------------------------------------------------------------------------
$ws.Range("G7").Select() | Out-Null
